# Applies the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list on Sun Sep  3 19:24:46 UTC 2023 with GitHub Actions".
# Column D = Price, Column E = Volume(1h), both stored as text in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row number -> @(new Price text (or $null to leave as-is), new Volume(1h) text (or $null))
$updates = @{
    2 = @('25.965.77', '  +0.40%  ')
    3 = @('1.640.86', '  +0.32%  ')
    4 = @($null, '  -0.26%  ')
    5 = @('214.69', '  +0.36%  ')
    6 = @('0.5091', '  +1.38%  ')
    7 = @($null, '  +0.05%  ')
    8 = @('0.2562', '  +0.11%  ')
    9 = @($null, '  +0.22%  ')
    10 = @('19.49', '  +0.31%  ')
    11 = @('0.07739', '  -0.49%  ')
    12 = @('4.278', '  +0.82%  ')
    13 = @('1.644.58', '  +0.60%  ')
    14 = @($null, '  +0.83%  ')
    15 = @('0.0₅7753', '  -1.36%  ')
    16 = @('64.26', '  -0.07%  ')
    17 = @('25.979.06', '  +0.41%  ')
    18 = @($null, '  -0.12%  ')
    19 = @('195.82', '  -0.06%  ')
    20 = @('4.429', '  +1.71%  ')
    21 = @('9.922', '  +0.44%  ')
    22 = @($null, '  +1.56%  ')
    23 = @('1.005', '  +0.40%  ')
    24 = @('1.886', '  +0.24%  ')
    25 = @('141.15', '  +1.30%  ')
    26 = @('0.1198', '  +5.82%  ')
    27 = @('6.846', '  +0.61%  ')
    28 = @('15.55', '  -0.61%  ')
    30 = @('0.04863', '  +0.45%  ')
    31 = @('3.250', '  +0.33%  ')
    32 = @($null, '  +0.30%  ')
    33 = @($null, '  +0.13%  ')
    34 = @('2.367', '  +0.67%  ')
    35 = @('0.8936', '  +1.05%  ')
    36 = @('1.144.17', '  +2.00%  ')
    37 = @('2.578', '  -0.86%  ')
    38 = @('0.5445', '  -1.03%  ')
    39 = @('0.01556', '  +0.21%  ')
    40 = @($null, '  +0.32%  ')
    41 = @('2.523', '  -1.69%  ')
    42 = @('0.0₈126', '  +4.03%  ')
    43 = @('0.8115', '  +0.23%  ')
    44 = @('99.06', $null)
    45 = @('5.436', '  -3.69%  ')
    46 = @('1.778.59', '  +0.36%  ')
    47 = @('0.4525', '  +0.12%  ')
    48 = @('54.89', '  -0.01%  ')
    49 = @('0.9970', '  -0.65%  ')
    50 = @('0.05053', '  +0.22%  ')
    51 = @('1.003', '  +0.01%  ')
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $dVal = $pair[0]
    $eVal = $pair[1]

    if ($dVal -ne $null) {
        # Force text format so price strings like "214.69" or "25.965.77"
        # are preserved verbatim instead of being parsed as numbers/dates.
        $dCell = $ws.Range("D$row")
        $dCell.NumberFormat = "@"
        $dCell.Value = $dVal
    }

    if ($eVal -ne $null) {
        $eCell = $ws.Range("E$row")
        $eCell.NumberFormat = "@"
        $eCell.Value = $eVal
    }
}
